# Automatische test-sync: 2025-06-19 22:35:50
# Adds the new incoming-mail log row (row 48) to the "Logs" sheet and
# updates the "Dashboard" category-count summary to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append the new log entry on the "Logs" sheet
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A48").Value = "Probleem met inloggen"
$logs.Range("B48").Value = "mailmind.test@zohomail.eu"
$logs.Range("C48").Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$logs.Range("D48").Value = "IT / Technisch probleem"
$logs.Range("F48").Value = "2025-06-19 22:35:14"
$logs.Range("G48").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too
# (Category column D and Answered column G), preserving existing rules.
$dRules = $logs.Range("D2:D47").FormatConditions
for ($i = 1; $i -le $dRules.Count; $i++) {
    $dRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D48"))
}

$gRules = $logs.Range("G2:G47").FormatConditions
for ($i = 1; $i -le $gRules.Count; $i++) {
    $gRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G48"))
}

# ---------------------------------------------------------------------
# 2. Update the "Dashboard" summary count for "IT / Technisch probleem"
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B8").Value = 4
